# BOM update: the MCU (row 3 - "ARM Microcontrollers - MCU 32KB Flash 4KB SRAM
# PSoC 4") unit price changed from $4.33 to $7.78 after the schematic gained
# the MCU, crystal and programming header. Updating the Unit Price cell lets
# the workbook's own formulas (Extended Price = Unit Price * Quantity, and the
# SUBTOTAL = SUM of all Extended Price cells) recompute on their own.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I3").Value = 7.78
